# Generate Report for Handoff
#
# The localization report previously listed "6c6133ff..." in row 5 and
# "ff7f135a..." in row 6 on every sheet, both shown as "In Translation".
# A handoff just happened for "6c6133ff...": it now needs to be reported as
# "Ready for handoff" with a fresh handoff timestamp, which also reorders it
# below "ff7f135a..." (now back to "In Translation") in the report - so the
# two rows swap places.
#
# NOTE: only the hyperlink *display text* moves with the row; the underlying
# hyperlink target (r:id / Address) stays attached to its original cell
# position - matching the source data exactly.

function Set-LinkText($ws, $row, $col, $text) {
    foreach ($h in $ws.Hyperlinks) {
        if (($h.Range.Row -eq $row) -and ($h.Range.Column -eq $col)) {
            $h.TextToDisplay = $text
        }
    }
}

$wb = $excel.ActiveWorkbook

$ff7f135a = "ff7f135a-bb2f-4e7b-8b06-4236cca36a5a.md"
$sixC = "6c6133ff-db15-4758-aab1-873abf3f3715.md"
$inTranslation = "In Translation"
$readyForHandoff = "Ready for handoff"

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Cells.Item(5, 1).Value = $ff7f135a
$wsOverview.Cells.Item(5, 2).Value = $inTranslation
$wsOverview.Cells.Item(5, 3).Value = $inTranslation

$wsOverview.Cells.Item(6, 1).Value = $sixC
$wsOverview.Cells.Item(6, 2).Value = $readyForHandoff
$wsOverview.Cells.Item(6, 3).Value = $readyForHandoff

Set-LinkText $wsOverview 5 1 $ff7f135a
Set-LinkText $wsOverview 6 1 $sixC

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Cells.Item(5, 1).Value = $ff7f135a
$wsZhCn.Cells.Item(5, 2).Value = $inTranslation
$wsZhCn.Cells.Item(5, 3).Value = "ff7f135a-bb2f-4e7b-8b06-4236cca36a5a.3e14b8dd899da539dd170c68c834efb820e8f44c.zh-cn.xlf"
$wsZhCn.Cells.Item(5, 4).Value = "2016-02-24 09:25:45"

$wsZhCn.Cells.Item(6, 1).Value = $sixC
$wsZhCn.Cells.Item(6, 2).Value = $readyForHandoff
$wsZhCn.Cells.Item(6, 3).Value = "6c6133ff-db15-4758-aab1-873abf3f3715.a28ad2c323baa1c0439506b7043c682751b6ccf5.zh-cn.xlf"
$wsZhCn.Cells.Item(6, 4).Value = "2016-02-24 09:35:38"

Set-LinkText $wsZhCn 5 1 $ff7f135a
Set-LinkText $wsZhCn 5 3 "ff7f135a-bb2f-4e7b-8b06-4236cca36a5a.3e14b8dd899da539dd170c68c834efb820e8f44c.zh-cn.xlf"
Set-LinkText $wsZhCn 6 1 $sixC
Set-LinkText $wsZhCn 6 3 "6c6133ff-db15-4758-aab1-873abf3f3715.a28ad2c323baa1c0439506b7043c682751b6ccf5.zh-cn.xlf"

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Cells.Item(5, 1).Value = $ff7f135a
$wsDeDe.Cells.Item(5, 2).Value = $inTranslation
$wsDeDe.Cells.Item(5, 3).Value = "ff7f135a-bb2f-4e7b-8b06-4236cca36a5a.3e14b8dd899da539dd170c68c834efb820e8f44c.de-de.xlf"
$wsDeDe.Cells.Item(5, 4).Value = "2016-02-24 09:25:57"

$wsDeDe.Cells.Item(6, 1).Value = $sixC
$wsDeDe.Cells.Item(6, 2).Value = $readyForHandoff
$wsDeDe.Cells.Item(6, 3).Value = "6c6133ff-db15-4758-aab1-873abf3f3715.a28ad2c323baa1c0439506b7043c682751b6ccf5.de-de.xlf"
$wsDeDe.Cells.Item(6, 4).Value = "2016-02-24 09:35:52"

Set-LinkText $wsDeDe 5 1 $ff7f135a
Set-LinkText $wsDeDe 5 3 "ff7f135a-bb2f-4e7b-8b06-4236cca36a5a.3e14b8dd899da539dd170c68c834efb820e8f44c.de-de.xlf"
Set-LinkText $wsDeDe 6 1 $sixC
Set-LinkText $wsDeDe 6 3 "6c6133ff-db15-4758-aab1-873abf3f3715.a28ad2c323baa1c0439506b7043c682751b6ccf5.de-de.xlf"
